$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the input parameters for the projection model
$ws.Range("E1").Value = 200    # point per month
$ws.Range("E3").Value = 2022   # current year
$ws.Range("E4").Value = 1      # current month
$ws.Range("E5").Value = 47000  # starting asset

# Update the selected cell to reflect where the author left off
$ws.Range("G25").Select()
